$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 10,23
$data[0,0] = 1.850998401641846
$data[0,1] = 3
$data[0,2] = 6108.976187777497
$data[0,3] = 0.2100780598911265
$data[0,4] = 0.1600597640567812
$data[0,5] = 0.1568647739662924
$data[0,6] = 0.155947264311399
$data[0,7] = 0.155829573518702
$data[0,8] = 0.1524256689263317
$data[0,9] = 0.1524256689263317
$data[0,10] = 0.1524224639274797
$data[0,11] = 0.1509460288757003
$data[0,12] = 0.1499289225666692
$data[0,13] = 0.1464880652105072
$data[0,14] = 0.1464880652105072
$data[0,15] = 0.1464880652105072
$data[0,16] = 0.1462150297469437
$data[0,17] = 0.1447380139405461
$data[0,18] = 0.1445939337924105
$data[0,19] = 0.1434838983807895
$data[0,20] = 0.1434838983807895
$data[0,21] = 0.1434838983807895
$data[0,22] = 0.1430833564868907
$data[1,0] = 1.959998369216919
$data[1,1] = 3
$data[1,2] = 6061.061698657956
$data[1,3] = 0.2100780598911265
$data[1,4] = 0.1600597640567812
$data[1,5] = 0.1568647739662924
$data[1,6] = 0.1540551784051227
$data[1,7] = 0.1540551784051227
$data[1,8] = 0.1540551784051227
$data[1,9] = 0.1540551784051227
$data[1,10] = 0.1540551784051227
$data[1,11] = 0.1540551784051227
$data[1,12] = 0.1533941637131498
$data[1,13] = 0.1483903279974793
$data[1,14] = 0.1483903279974793
$data[1,15] = 0.1450240524904412
$data[1,16] = 0.1450240524904412
$data[1,17] = 0.1443688032806061
$data[1,18] = 0.1443688032806061
$data[1,19] = 0.143496164911475
$data[1,20] = 0.1427066977472081
$data[1,21] = 0.1425745785387607
$data[1,22] = 0.1421493508510323
$data[2,0] = 1.929001569747925
$data[2,1] = 3
$data[2,2] = 6064.721418270442
$data[2,3] = 0.2100780598911265
$data[2,4] = 0.1600597640567812
$data[2,5] = 0.1568647739662924
$data[2,6] = 0.1560440163731134
$data[2,7] = 0.151221569005294
$data[2,8] = 0.147242221359071
$data[2,9] = 0.1429953943945027
$data[2,10] = 0.1429953943945027
$data[2,11] = 0.1429953943945027
$data[2,12] = 0.1429953943945027
$data[2,13] = 0.1429953943945027
$data[2,14] = 0.1429953943945027
$data[2,15] = 0.1429953943945027
$data[2,16] = 0.1429953943945027
$data[2,17] = 0.1429953943945027
$data[2,18] = 0.1429953943945027
$data[2,19] = 0.1422354836039892
$data[2,20] = 0.1422354836039892
$data[2,21] = 0.1422354836039892
$data[2,22] = 0.1422206904146285
$data[3,0] = 1.777008533477783
$data[3,1] = 3
$data[3,2] = 6094.222036279522
$data[3,3] = 0.2100780598911265
$data[3,4] = 0.1600597640567812
$data[3,5] = 0.1568647739662924
$data[3,6] = 0.1560440163731134
$data[3,7] = 0.1537640374119212
$data[3,8] = 0.1537640374119212
$data[3,9] = 0.1537640374119212
$data[3,10] = 0.1537640374119212
$data[3,11] = 0.1509546934121953
$data[3,12] = 0.1476937255697986
$data[3,13] = 0.1476937255697986
$data[3,14] = 0.1476937255697986
$data[3,15] = 0.1452758186161132
$data[3,16] = 0.1447366241578489
$data[3,17] = 0.1435533082807948
$data[3,18] = 0.1435533082807948
$data[3,19] = 0.1435533082807948
$data[3,20] = 0.1435533082807948
$data[3,21] = 0.1435533082807948
$data[3,22] = 0.1427957511945326
$data[4,0] = 1.849978685379028
$data[4,1] = 3
$data[4,2] = 6207.786417438995
$data[4,3] = 0.2100780598911265
$data[4,4] = 0.1600597640567812
$data[4,5] = 0.1568647739662924
$data[4,6] = 0.1560440163731134
$data[4,7] = 0.147381492381348
$data[4,8] = 0.147381492381348
$data[4,9] = 0.147381492381348
$data[4,10] = 0.147381492381348
$data[4,11] = 0.1464383651948443
$data[4,12] = 0.1462087375914113
$data[4,13] = 0.1462087375914113
$data[4,14] = 0.1462087375914113
$data[4,15] = 0.1462087375914113
$data[4,16] = 0.1462087375914113
$data[4,17] = 0.1458166668098168
$data[4,18] = 0.1458166668098168
$data[4,19] = 0.1456608945546251
$data[4,20] = 0.1456608945546251
$data[4,21] = 0.1451553182299791
$data[4,22] = 0.1450094818214229
$data[5,0] = 1.672035694122314
$data[5,1] = 3
$data[5,2] = 6118.56575990176
$data[5,3] = 0.2100780598911265
$data[5,4] = 0.1600597640567812
$data[5,5] = 0.1568647739662924
$data[5,6] = 0.1560440163731134
$data[5,7] = 0.1551958969359106
$data[5,8] = 0.1541769242807673
$data[5,9] = 0.1541769242807673
$data[5,10] = 0.1480694999173292
$data[5,11] = 0.1471405351641092
$data[5,12] = 0.1462608025093261
$data[5,13] = 0.1462608025093261
$data[5,14] = 0.1445759192040776
$data[5,15] = 0.1434925387287085
$data[5,16] = 0.1434925387287085
$data[5,17] = 0.1434925387287085
$data[5,18] = 0.1432702877173832
$data[5,19] = 0.1432702877173832
$data[5,20] = 0.1432702877173832
$data[5,21] = 0.1432702877173832
$data[5,22] = 0.1432702877173832
$data[6,0] = 1.82099461555481
$data[6,1] = 3
$data[6,2] = 6173.329300477842
$data[6,3] = 0.2100780598911265
$data[6,4] = 0.1600597640567812
$data[6,5] = 0.1568647739662924
$data[6,6] = 0.1560440163731134
$data[6,7] = 0.155829573518702
$data[6,8] = 0.1557735451149562
$data[6,9] = 0.1549212787294311
$data[6,10] = 0.1549212787294311
$data[6,11] = 0.1510225897215471
$data[6,12] = 0.1481328739891422
$data[6,13] = 0.1481328739891422
$data[6,14] = 0.1481328739891422
$data[6,15] = 0.1478390686684193
$data[6,16] = 0.1478390686684193
$data[6,17] = 0.146347800249354
$data[6,18] = 0.1452363234564854
$data[6,19] = 0.1452363234564854
$data[6,20] = 0.1452363234564854
$data[6,21] = 0.144337803128223
$data[6,22] = 0.144337803128223
$data[7,0] = 1.636015176773071
$data[7,1] = 3
$data[7,2] = 5963.271113327468
$data[7,3] = 0.2100780598911265
$data[7,4] = 0.1600597640567812
$data[7,5] = 0.1568647739662924
$data[7,6] = 0.1560440163731134
$data[7,7] = 0.1549574303356408
$data[7,8] = 0.1488304963110333
$data[7,9] = 0.1441655968877454
$data[7,10] = 0.1441655968877454
$data[7,11] = 0.1441655968877454
$data[7,12] = 0.1441655968877454
$data[7,13] = 0.1427722349072908
$data[7,14] = 0.1427722349072908
$data[7,15] = 0.1427722349072908
$data[7,16] = 0.1419677473112088
$data[7,17] = 0.1411242855636162
$data[7,18] = 0.1411242855636162
$data[7,19] = 0.1406520023074759
$data[7,20] = 0.1405939555667922
$data[7,21] = 0.1404101780994849
$data[7,22] = 0.1402431016243171
$data[8,0] = 1.819993495941162
$data[8,1] = 3
$data[8,2] = 6181.991474539906
$data[8,3] = 0.2100780598911265
$data[8,4] = 0.1600597640567812
$data[8,5] = 0.1568647739662924
$data[8,6] = 0.1560440163731134
$data[8,7] = 0.145129548347398
$data[8,8] = 0.145129548347398
$data[8,9] = 0.145129548347398
$data[8,10] = 0.145129548347398
$data[8,11] = 0.145129548347398
$data[8,12] = 0.145129548347398
$data[8,13] = 0.145129548347398
$data[8,14] = 0.145129548347398
$data[8,15] = 0.145129548347398
$data[8,16] = 0.145129548347398
$data[8,17] = 0.145129548347398
$data[8,18] = 0.1445066564237798
$data[8,19] = 0.1445066564237798
$data[8,20] = 0.1445066564237798
$data[8,21] = 0.1445066564237798
$data[8,22] = 0.1445066564237798
$data[9,0] = 1.697025299072266
$data[9,1] = 3
$data[9,2] = 6104.453648436408
$data[9,3] = 0.2100780598911265
$data[9,4] = 0.1600597640567812
$data[9,5] = 0.1533160344061542
$data[9,6] = 0.1533160344061542
$data[9,7] = 0.1532006730738767
$data[9,8] = 0.1506392647723403
$data[9,9] = 0.1434340953178387
$data[9,10] = 0.1434340953178387
$data[9,11] = 0.1434340953178387
$data[9,12] = 0.1434340953178387
$data[9,13] = 0.1434340953178387
$data[9,14] = 0.1434340953178387
$data[9,15] = 0.1434340953178387
$data[9,16] = 0.1434340953178387
$data[9,17] = 0.1433527779795554
$data[9,18] = 0.1433527779795554
$data[9,19] = 0.1433527779795554
$data[9,20] = 0.1429951978252711
$data[9,21] = 0.1429951978252711
$data[9,22] = 0.1429951978252711

$ws.Range("C2:Y11").Value = $data

$wb.Save()
